$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("C10").Value = 345539
$ws.Range("D10").Value = 64165
$ws.Range("E10").Value = 1817720252

# Row 21
$ws.Range("C21").Value = 175240
$ws.Range("D21").Value = 38057
$ws.Range("E21").Value = 316813398

# Row 78
$ws.Range("C78").Value = 178441
$ws.Range("E78").Value = 892526355

# Row 121
$ws.Range("C121").Value = 1306164
$ws.Range("E121").Value = 2274669692

# Row 129
$ws.Range("C129").Value = 633440
$ws.Range("E129").Value = 3428621480

# Row 132
$ws.Range("C132").Value = 585716
$ws.Range("E132").Value = 3463679198

# Row 144
$ws.Range("C144").Value = 25079
$ws.Range("E144").Value = 92448942

# Row 154
$ws.Range("C154").Value = 18455
$ws.Range("E154").Value = 72768683

# Row 156
$ws.Range("C156").Value = 12402
$ws.Range("E156").Value = 40299139

# Row 194
$ws.Range("C194").Value = 18380
$ws.Range("E194").Value = 71349778

# Row 229
$ws.Range("C229").Value = 612546
$ws.Range("E229").Value = 1040763491
